$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 881
$ws.Range("F3").Value = 1821
$ws.Range("F4").Value = 100
$ws.Range("F5").Value = 550
$ws.Range("F6").Value = 2147
$ws.Range("F8").Value = 2126
$ws.Range("F10").Value = 55
$ws.Range("F11").Value = 2455
$ws.Range("F12").Value = 674
$ws.Range("F14").Value = 4138
$ws.Range("F17").Value = 3251
$ws.Range("F18").Value = 894
$ws.Range("F19").Value = 154
$ws.Range("F21").Value = 216
$ws.Range("F22").Value = 2195
$ws.Range("F23").Value = 1219
$ws.Range("F24").Value = 11
$ws.Range("F25").Value = 1998
$ws.Range("F26").Value = 438
$ws.Range("F29").Value = 9070
$ws.Range("F30").Value = 5832
$ws.Range("F31").Value = 367
$ws.Range("F34").Value = 46
$ws.Range("F35").Value = 792
$ws.Range("F38").Value = 950
$ws.Range("F39").Value = 417
$ws.Range("F40").Value = 77
$ws.Range("F41").Value = 211
$ws.Range("F43").Value = 4708
$ws.Range("F44").Value = 7
$ws.Range("F45").Value = 918
$ws.Range("F46").Value = 100

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("G2").Value = 224
$ws.Range("F12").Value = 84
$ws.Range("F13").Value = 115
$ws.Range("F16").Value = 3478
$ws.Range("F23").Value = 8

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 8503
$ws.Range("F3").Value = 384
$ws.Range("F4").Value = 1410

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 8503
$ws.Range("F3").Value = 881
$ws.Range("F4").Value = 384
$ws.Range("F5").Value = 1410
$ws.Range("F6").Value = 1821
$ws.Range("F7").Value = 100
$ws.Range("F9").Value = 55
$ws.Range("F11").Value = 4138
$ws.Range("F14").Value = 3251
$ws.Range("F15").Value = 894
$ws.Range("F16").Value = 154
$ws.Range("F17").Value = 216
$ws.Range("F18").Value = 2195
$ws.Range("F23").Value = 1219
$ws.Range("F24").Value = 84
$ws.Range("F25").Value = 11
$ws.Range("F26").Value = 115
$ws.Range("F27").Value = 438
$ws.Range("F30").Value = 9070
$ws.Range("F31").Value = 3478
$ws.Range("F33").Value = 367
$ws.Range("F35").Value = 792
$ws.Range("F36").Value = 950
$ws.Range("F37").Value = 417
$ws.Range("F38").Value = 77
$ws.Range("F39").Value = 211
$ws.Range("F42").Value = 4708
$ws.Range("F43").Value = 7
$ws.Range("F44").Value = 918
$ws.Range("F47").Value = 8
